$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "41.933.97"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.287.12"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "316.44"
$ws.Range("E5").Value = "  -0.28%  "
Set-TextValue $ws.Range("D6") "103.33"
$ws.Range("E6").Value = "  -3.33%  "
Set-TextValue $ws.Range("D7") "0.624"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue $ws.Range("D9") "0.603"
$ws.Range("E9").Value = "  -2.54%  "
Set-TextValue $ws.Range("D10") "39.25"
$ws.Range("E10").Value = "  -4.76%  "
Set-TextValue $ws.Range("D11") "0.0903"
$ws.Range("E11").Value = "  -2.50%  "
Set-TextValue $ws.Range("D12") "8.23"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("E13").Value = "  -0.80%  "
Set-TextValue $ws.Range("D14") "0.957"
$ws.Range("E14").Value = "  -4.27%  "
Set-TextValue $ws.Range("D15") "15.16"
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("D16").Value = "2.630.98"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "2.285.64"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "42.012.68"
$ws.Range("E18").Value = "  -1.87%  "
Set-TextValue $ws.Range("D19") "7.36"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("E20").Value = "  -1.07%  "
Set-TextValue $ws.Range("D21") "3.61"
$ws.Range("E21").Value = "  -0.43%  "
Set-TextValue $ws.Range("D22") "73.13"
$ws.Range("E22").Value = "  -3.79%  "
Set-TextValue $ws.Range("D23") "277.89"
$ws.Range("E23").Value = "  +3.30%  "
Set-TextValue $ws.Range("D24") "10.10"
$ws.Range("E24").Value = "  +7.78%  "
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("E28").Value = "  +4.49%  "
Set-TextValue $ws.Range("D29") "22.78"
$ws.Range("E29").Value = "  -2.86%  "
Set-TextValue $ws.Range("D30") "36.03"
$ws.Range("E30").Value = "  -1.77%  "
Set-TextValue $ws.Range("D31") "162.85"
$ws.Range("E31").Value = "  -2.66%  "
Set-TextValue $ws.Range("D32") "0.0868"
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("E35").Value = "  +4.11%  "
Set-TextValue $ws.Range("D36") "0.113"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("E39").Value = "  +4.34%  "
Set-TextValue $ws.Range("D40") "3.74"
$ws.Range("E40").Value = "  -2.09%  "
Set-TextValue $ws.Range("D41") "99.38"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("E42").Value = "  -4.55%  "
Set-TextValue $ws.Range("D43") "69.06"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -6.11%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "112.55"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D47") "11.87"
$ws.Range("E47").Value = "  -3.68%  "
Set-TextValue $ws.Range("D48") "77.02"
$ws.Range("E48").Value = "  +1.95%  "
Set-TextValue $ws.Range("D49") "8.87"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("E50").Value = "  -4.67%  "
$ws.Range("D51").Value = "1.575.09"
$ws.Range("E51").Value = "  +0.09%  "
